# Add a new "AllowFields" row to the Configuration sheet, under the
# existing XNAT config block, mirroring the existing Project/Server/Scans
# rows (bold label in column B, plain value in column C).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Configuration")

$ws.Range("B13").Value = "AllowFields"
$ws.Range("C13").Value = "AccessionNumber"

# Match the bold styling used by the other config-label cells in column B
$ws.Range("B13").Font.Bold = $true

# Leave the selection where the user ended up after typing the new row
$ws.Range("B13").Select()
